$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Backend and Frontend use different account transaction format"
$ws.Range("B3").Value = "The Backend uses spaces, and the Frontend uses underscores"
$ws.Range("C3").Value = "Updated the Backend to read in input files with underscores"

$ws.Range("A4").Value = "Backend expects login transaction?"

$ws.Range("A5").Value = "Frontend has END_OF_FILE line in user accounts"
$ws.Range("B5").Value = "Backend not expecting 00000_END_OF_FILE____________00000.00"
$ws.Range("C5").Value = "Provide adapter to remove line"

$ws.Range("A6").Value = "FrontEnd file output does match Backend file input"
$ws.Range("B6").Value = "Front end user accounts doesn't match the back end expected input for user accounts"

$ws.Columns.Item(2).ColumnWidth = 60.8
$ws.Range("B16").Select()
